$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> (DAMSLTag, DialogAct)
$updates = @{
    15  = @("sd", "Statement-non-opinion")
    26  = @("sd", "Statement-non-opinion")
    45  = @("%", "Uninterpretable")
    47  = @("%", "Uninterpretable")
    57  = @("sd", "Statement-non-opinion")
    63  = @("sd", "Statement-non-opinion")
    70  = @("aa", "Agree/Accept")
    72  = @("sv", "Statement-opinion")
    74  = @("b", "Acknowledge (Backchannel)")
    75  = @("sd", "Statement-non-opinion")
    82  = @("sv", "Statement-opinion")
    102 = @("sd", "Statement-non-opinion")
    105 = @("sd", "Statement-non-opinion")
    111 = @("b", "Acknowledge (Backchannel)")
    113 = @("sd", "Statement-non-opinion")
    124 = @("sv", "Statement-opinion")
    157 = @("sv", "Statement-opinion")
    159 = @("sd", "Statement-non-opinion")
}

foreach ($row in $updates.Keys) {
    $values = $updates[$row]
    $ws.Range("I$row").Value = $values[0]
    $ws.Range("J$row").Value = $values[1]
}

$wb.Save()
